$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1!A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.66 = 13936.47 pesos`n✅ 13936.47 pesos = 3.65 = 944.52 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 273.24
$ws2.Range("O10").Value = 3808
$ws2.Range("N12").Value = 3822.9
$ws2.Range("O12").Value = 259.09
